$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the multi-word column headers to use underscores instead of spaces.
# (The shared-string table is rebuilt by the engine on save, so renaming the
# text here reproduces the upstream reorder/append of the shared strings.)
$ws.Range("F1").Value = "AVG_POSS"
$ws.Range("G1").Value = "TOU_IN_BOX"
$ws.Range("I1").Value = "FAST_BREAKS"
$ws.Range("J1").Value = "PASS_F3RD_TOT"
$ws.Range("K1").Value = "PASS_F3RD_SUC"
$ws.Range("L1").Value = "CROSSES_TOT"
$ws.Range("M1").Value = "CROSS_SUC"
$ws.Range("N1").Value = "TROUGH_BALL"

# Move the active selection to N1, matching the saved view state.
$ws.Range("N1").Select()
